# Applies the "case with 380 kV done" update to Case_4_46 res_bus/vm_pu.xlsx.
#
# The power-flow slack/reference voltage setpoint (column B, bus "1") was
# changed from 1.05 pu to 1.02 pu, and the case was re-run for the 380 kV
# network. That changes every downstream bus voltage magnitude, so columns
# B, C, D, E, F, I, J, K, L, M, N (rows 2-25, i.e. bus indices 0-23) are
# rewritten with the freshly computed results. Column G (constant 1, a
# reference bus) and the header row/column A (bus index) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1-based sheet column indices for B,C,D,E,F,I,J,K,L,M,N (G/H are skipped -
# they are not part of this edit).
$colIndex = @(2,3,4,5,6,9,10,11,12,13,14)

# Each row: @(rowNumber, valueForB, valueForC, ..., valueForN)
$data = @(
    @(2,1.02,1.086111515195323,1.079326327007237,1.099526299028462,1.105858638932316,1.067993881393705,1.090963080620817,1.082002670846362,1.102150772160547,1.108467271695202,1.092512373330512),
    @(3,1.02,1.08767390701954,1.080566316421894,1.101077789504023,1.107469600343693,1.068629136575458,1.092186123496859,1.083059683237019,1.103522324142208,1.109899314204997,1.093737153067724),
    @(4,1.02,1.088681789224686,1.081365366605998,1.1020789707177,1.108509433274302,1.06903707286139,1.092973980237557,1.08373977272769,1.10440656443865,1.110822874699404,1.094526128655362),
    @(5,1.02,1.089104775717904,1.081700505934934,1.10249922390309,1.108945976612105,1.069207829364179,1.093304359584229,1.084024765864963,1.104777533147244,1.111210416650344,1.094856977178601),
    @(6,1.02,1.089175754730363,1.081756731773883,1.102569748874368,1.10901923918176,1.069236456920753,1.093359783028166,1.084072564000112,1.104839775832349,1.111275444585799,1.094912479330186),
    @(7,1.02,1.088687444031045,1.081369847810485,1.102084588673171,1.10851526873824,1.069039357417028,1.092978398051729,1.083743584407922,1.104411524334948,1.110828055878775,1.094530552743337),
    @(8,1.02,1.086640180151288,1.079746077827358,1.100051206328459,1.10640360966252,1.068209217816452,1.091377151648578,1.082360699858036,1.102614973284847,1.108951878998046,1.092927032386627),
    @(9,1.02,1.08300841135702,1.076859014312974,1.09644663283486,1.102662393154863,1.066722261739977,1.088528018635396,1.079893822889963,1.099423864036657,1.105621817744058,1.090073853277768),
    @(10,1.02,1.080570149147531,1.074916337506234,1.094028383052663,1.100153899285964,1.065714334689234,1.086609427482301,1.07822842149857,1.097278706459719,1.103384924575015,1.08815253750517),
    @(11,1.02,1.079510125078636,1.074070731846755,1.092977480771924,1.099064118556567,1.065273860851421,1.085773962095441,1.077502210554354,1.09634546152598,1.102412167112976,1.087315885661631),
    @(12,1.02,1.079115732434404,1.073755960657315,1.092586545756304,1.09865877176214,1.065109635445442,1.085462913684503,1.077231688229864,1.095998141430697,1.102050201514255,1.087004395526271),
    @(13,1.02,1.079200360772319,1.073823510920015,1.092670429267512,1.09874574532313,1.065144890252196,1.085529667384998,1.077289751460133,1.096072673372834,1.102127873558847,1.087071244024678),
    @(14,1.02,1.079477537860463,1.074044726616702,1.092945177941979,1.099030623846405,1.065260298494144,1.085748265486385,1.077479864978481,1.096316765714181,1.102382260064271,1.087290152560442),
    @(15,1.02,1.079648228884333,1.074180935009611,1.093114381960202,1.099206073086196,1.065331323755286,1.085882855281785,1.077596897107025,1.096467069685705,1.102538910737098,1.08742493348879),
    @(16,1.02,1.080640408799675,1.074972363514012,1.094098047006491,1.100226147615791,1.065743481829792,1.086664774394971,1.07827650966402,1.097340549456971,1.103449394280645,1.088207963016803),
    @(17,1.02,1.081261630923019,1.075467615248368,1.094714050725826,1.100865042787386,1.06600093211773,1.087153983522927,1.07870144396746,1.097887278626128,1.104019390854982,1.088697866877834),
    @(18,1.02,1.081623571187335,1.075756062190344,1.095072991150223,1.101237354518349,1.066150709737984,1.087438878057048,1.078948811508049,1.098205755064141,1.104351459018068,1.088983165994872),
    @(19,1.02,1.081746914709075,1.075854343428589,1.09519531923629,1.101364245125223,1.066201714322065,1.087535943301576,1.079033074784457,1.09831427622838,1.104464618050992,1.08908036908318),
    @(20,1.02,1.081195021990425,1.075414523475726,1.094647997124842,1.100796531175879,1.065973350406439,1.087101542948784,1.078655903249889,1.097828663450508,1.103958277215357,1.088645351832062),
    @(21,1.02,1.079395934253642,1.073979602858384,1.092864287520159,1.098946749627326,1.065226330649271,1.085683913755777,1.077423902782271,1.096244905247493,1.102307367363579,1.087225709442996),
    @(22,1.02,1.078260993411612,1.073073496233953,1.091739417413245,1.097780507454907,1.064753094515765,1.084788425517209,1.076644803886797,1.095245242343522,1.10126566397704,1.086328949508415),
    @(23,1.02,1.07886301080439,1.073554215806679,1.092336057741369,1.098399063890686,1.065004305503197,1.085263540223106,1.077058248670659,1.095775556023331,1.101818247104029,1.086804738931691),
    @(24,1.02,1.081225120960746,1.075438514675652,1.0946778450289,1.100827489690043,1.065985814595747,1.087125239986191,1.078676482652123,1.097855150427409,1.103985893078057,1.08866908252198),
    @(25,1.02,1.083950262960047,1.077608509076448,1.097381125284486,1.103632059843853,1.067109576588241,1.089267917819557,1.08053519384912,1.100251919066999,1.106485635654225,1.090814803203636)
)

foreach ($entry in $data) {
    $row = $entry[0]
    for ($i = 0; $i -lt $colIndex.Length; $i++) {
        $col = $colIndex[$i]
        $value = $entry[$i + 1]
        $ws.Cells.Item($row, $col).Value = $value
    }
}
